# Refresh the crypto symbol/price snapshot (GitHub Actions scheduled update).
# Price cells in column D are stored as text (e.g. "0.00005920" keeps trailing
# zeros that a numeric value would drop), so each new price is written with a
# leading apostrophe to force Excel to keep it as text instead of a number.
# A couple of rows also got their column-E label text corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.11"
$ws.Range("D3").Value = "'21.77"
$ws.Range("D5").Value = "'0.05701"
$ws.Range("D6").Value = "'3.382"
$ws.Range("D7").Value = "'0.8099"
$ws.Range("D8").Value = "'1.023"
$ws.Range("D9").Value = "'0.1454"
$ws.Range("D10").Value = "'0.07509"
$ws.Range("D11").Value = "'0.03178"
$ws.Range("D12").Value = "'0.03017"
$ws.Range("D13").Value = "'0.09273"
$ws.Range("D14").Value = "'3.619"
$ws.Range("D15").Value = "'0.001661"
$ws.Range("D16").Value = "'0.04693"
$ws.Range("D17").Value = "'0.0005851"
$ws.Range("D18").Value = "'0.006360"
$ws.Range("D19").Value = "'0.005038"
$ws.Range("D20").Value = "'0.001041"
$ws.Range("D21").Value = "'0.0001499"
$ws.Range("D22").Value = "'0.0003101"
$ws.Range("E22").Value = "21UpBotsUBXTWorstin24h"
$ws.Range("D23").Value = "'3.771"
$ws.Range("D24").Value = "'6.413"
$ws.Range("D25").Value = "'2.161"
$ws.Range("D26").Value = "'0.3291"
$ws.Range("D40").Value = "'0.04064"
$ws.Range("D41").Value = "'0.006973"
$ws.Range("D43").Value = "'0.002908"
$ws.Range("D44").Value = "'0.008526"
$ws.Range("D45").Value = "'0.00005920"
$ws.Range("D47").Value = "'0.0005501"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("D48").Value = "'0.6826"
$ws.Range("D49").Value = "'0.007958"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D51").Value = "'0.01010"
